$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit cyclically rotates the data of rows 2-5 (row3->row2, row4->row3,
# row5->row4, row2->row5) and swaps the data of rows 11<->12 and 16<->17.
# Below we set every cell that actually changes value to its final
# (after-edit) content, and clear/add the "M" (Aktivitet) cell where the
# diff shows it disappearing/appearing.

# --- Row 2 (becomes old row 3's species, new Q/R) ---
$ws.Range("A2").Value = 130670771
$ws.Range("B2").Value = 79243
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("M2").ClearContents()
$ws.Range("Q2").Value = 491374
$ws.Range("R2").Value = 6759416

# --- Row 3 (becomes old row 4's id, new Q/R) ---
$ws.Range("A3").Value = 130670627
$ws.Range("Q3").Value = 491376
$ws.Range("R3").Value = 6759442

# --- Row 4 (becomes old row 5's species, new Q/R) ---
$ws.Range("A4").Value = 130661609
$ws.Range("B4").Value = 8451
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 106545
$ws.Range("F4").Value = "Mindre märgborre"
$ws.Range("G4").Value = "Tomicus minor"
$ws.Range("H4").Value = "(Hartig, 1834)"
$ws.Range("Q4").Value = 491477
$ws.Range("R4").Value = 6759416

# --- Row 5 (becomes old row 2's species, new Q/R) ---
$ws.Range("A5").Value = 130667093
$ws.Range("B5").Value = 57884
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("M5").Value = "färska spår"
$ws.Range("Q5").Value = 491408
$ws.Range("R5").Value = 6759381

# --- Rows 11 <-> 12 swap ---
$ws.Range("A11").Value = 130668644
$ws.Range("B11").Value = 79243
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = "Garnlav"
$ws.Range("G11").Value = "Alectoria sarmentosa"
$ws.Range("H11").Value = "(Ach.) Ach."
$ws.Range("M11").ClearContents()
$ws.Range("Q11").Value = 491424
$ws.Range("R11").Value = 6759256

$ws.Range("A12").Value = 130662234
$ws.Range("B12").Value = 8451
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 106545
$ws.Range("F12").Value = "Mindre märgborre"
$ws.Range("G12").Value = "Tomicus minor"
$ws.Range("H12").Value = "(Hartig, 1834)"
$ws.Range("M12").Value = "färska gnagspår"
$ws.Range("Q12").Value = 491455
$ws.Range("R12").Value = 6759425

# --- Rows 16 <-> 17 swap ---
$ws.Range("A16").Value = 130661510
$ws.Range("B16").Value = 79243
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 6425
$ws.Range("F16").Value = "Garnlav"
$ws.Range("G16").Value = "Alectoria sarmentosa"
$ws.Range("H16").Value = "(Ach.) Ach."
$ws.Range("M16").ClearContents()
$ws.Range("Q16").Value = 491504
$ws.Range("R16").Value = 6759336

$ws.Range("A17").Value = 130661548
$ws.Range("B17").Value = 8451
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 106545
$ws.Range("F17").Value = "Mindre märgborre"
$ws.Range("G17").Value = "Tomicus minor"
$ws.Range("H17").Value = "(Hartig, 1834)"
$ws.Range("M17").Value = "äldre gnagspår"
$ws.Range("Q17").Value = 491487
$ws.Range("R17").Value = 6759357
